$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 3
    18 = 0
    19 = 2
    20 = 2
    21 = 1
    22 = 2
    23 = 0
    24 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
